$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Z2:Z54").Value = "2025-11-05T14:13:28.231820"
$ws.Range("Z55").Value = "2025-11-05T14:13:28.330054"
$ws.Range("Z56:Z57").Value = "2025-11-05T14:13:28.331115"
$ws.Range("Z58:Z67").Value = "2025-11-05T14:13:28.331625"
